$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.895.47"
$ws.Range("E2").Value = "  +5.53%  "
$ws.Range("D3").Value = "3.103.86"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.78%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.098.47"
$ws.Range("E8").Value = "  +3.43%  "
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.85%  "
$ws.Range("E11").Value = "  +5.80%  "
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("E13").Value = "  +5.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").Value = "3.602.86"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "64.733.49"
$ws.Range("E16").Value = "  +5.23%  "
$ws.Range("D17").Value = "3.099.81"
$ws.Range("E17").Value = "  +3.29%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.676"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.96%  "
$ws.Range("E29").Value = "  +9.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("E35").Value = "  +6.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "464.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0408"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("D40").Value = "3.017.21"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.44%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.261"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.88%  "
$ws.Range("E48").Value = "  +4.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +7.73%  "
$ws.Range("E51").Value = "  +3.89%  "
